# Updated cryptos list - apply diff changes to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are kept as Text so values like "1.00" / "67.807.18"
# are not coerced into numbers by Excel's automatic type inference.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.807.18'
$ws.Range("E2").Value = '  -2.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.414.41'
$ws.Range("E3").Value = '  -2.98%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.54'
$ws.Range("E5").Value = '  -3.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.58'
$ws.Range("E6").Value = '  -3.24%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.495'
$ws.Range("E8").Value = '  -3.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.412.41'
$ws.Range("E9").Value = '  -3.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  -8.81%  '

$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.328'
$ws.Range("E12").Value = '  -6.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.72'
$ws.Range("E13").Value = '  -3.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.869.07'
$ws.Range("E14").Value = '  -2.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.718.60'
$ws.Range("E15").Value = '  -2.32%  '

$ws.Range("E16").Value = '  -5.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '22.91'
$ws.Range("E17").Value = '  -4.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.436.78'
$ws.Range("E18").Value = '  -2.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.53'
$ws.Range("E19").Value = '  -5.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.74'
$ws.Range("E20").Value = '  -3.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.85'
$ws.Range("E21").Value = '  -7.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.72'
$ws.Range("E22").Value = '  -4.04%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  -3.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.66'
$ws.Range("E25").Value = '  -5.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.572.04'
$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("E27").Value = '  -7.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.90'
$ws.Range("E29").Value = '  -8.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0790'
$ws.Range("E30").Value = '  -8.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.99'
$ws.Range("E31").Value = '  -7.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '426.10'
$ws.Range("E33").Value = '  -2.65%  '

$ws.Range("E34").Value = '  -7.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  -7.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.91'
$ws.Range("E36").Value = '  +0.32%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.98'
$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("E38").Value = '  -0.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.107'
$ws.Range("E39").Value = '  -5.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.54'
$ws.Range("E40").Value = '  -3.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.296'
$ws.Range("E41").Value = '  -5.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.27'
$ws.Range("E42").Value = '  -0.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.31'
$ws.Range("E43").Value = '  -5.52%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.42'
$ws.Range("E44").Value = '  -9.56%  '

$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.06'
$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.97'
$ws.Range("E46").Value = '  -8.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '129.15'
$ws.Range("E47").Value = '  -6.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.27'
$ws.Range("E48").Value = '  -4.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0706'
$ws.Range("E49").Value = '  -2.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.472'
$ws.Range("E50").Value = '  -6.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.548'
$ws.Range("E51").Value = '  -4.07%  '

Write-Host "Applied cryptos update."